# Map path fixes made during the meeting with Sean
#
# The source workbook's absolute-path bookmark moved from the author's
# "Projects" folder to their "Downloads" folder, the book window was
# repositioned/resized, and a new lookup row (Code "015" -> Value "F") was
# added below the existing table; the active selection and zoom level were
# also updated to reflect where the author was working when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data row -------------------------------------------------------
# Row 7: Code "015" maps to Value "F" (quotes are part of the literal text,
# matching the existing rows in columns A/B).
$ws.Range("A7").Value = '"015"'
$ws.Range("B7").Value = '"F"'

# --- View state ----------------------------------------------------------
# Zoom in on the sheet and leave the selection on B9, where the author's
# cursor ended up.
$win = $excel.ActiveWindow
$win.Zoom = 329

# Reposition/resize the workbook window (best-effort: reflects the window
# geometry recorded in the saved file).
$win.Left = -33920
$win.Top = -20240
$win.Width = 51200
$win.Height = 28300

$ws.Range("B9").Select()
